# Security of nature for all mun (08-22)
# Adds a new "Охрана окруж. среды" category (column D) to the existing
# category block in rows 24-25, matching the header/data style already
# used by the sibling B/C columns in that block, and updates the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New category header (row 24, column D) - same style as C24 (bold header)
$ws.Range("C24").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D24").Value = "Охрана окруж. среды"

# New category data row (row 25, column D) - same style as C25
$ws.Range("C25").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D25").Value = "Затраты на прир. охр. - naturesecure (тыс. руб.) (8055001)"

# Match the author's final selection
$ws.Range("E22").Select()
